{"js": "// \"fix for font size in document generation\"\n//\n// The \"\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \" label and the quoted \"$MAN_NAME\" placeholder used to\n// live in two separate runs, split apart by a leftover/empty DDE\n// bookmark. That split run boundary made the generator lose the run\n// formatting for part of the line. Collapse the two runs into a single\n// run reading  \u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: $MAN_NAME  (no quotes) and drop the\n// now-pointless bookmark that used to sit between them.\n\nconst body = context.document.body;\n\n// Locate the exact paragraph text we need to touch: '\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \"$MAN_NAME\"'\nconst results = body.search('\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \"$MAN_NAME\"', { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Replace the matched range with a single, self-contained run carrying\n  // the same run formatting (FreeSerif, 12pt/24 half-points, zxx lang)\n  // that the original two runs used - and no bookmark in between, so the\n  // OOXML serializer merges everything into one <w:r>.\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-16\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r><w:rPr>\" +\n    '<w:rFonts w:ascii=\"FreeSerif\" w:hAnsi=\"FreeSerif\"/>' +\n    '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n    '<w:lang w:val=\"zxx\" w:eastAsia=\"zxx\" w:bidi=\"zxx\"/>' +\n    \"</w:rPr><w:t>\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: $MAN_NAME</w:t></w:r></w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"fix for font size in document generation\"\n#\n# The \"\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \" label and the quoted \"$MAN_NAME\" placeholder used to\n# live in two separate runs (split apart by a leftover, empty DDE\n# bookmark), which made Word fall back to default/mismatched run\n# properties for part of the line when the template was regenerated.\n# Collapse them into a single run reading  \u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: $MAN_NAME\n# (no quotes) and drop the now-pointless bookmark in between.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the empty bookmark that sits between the \"\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \" run\n#    and the \"$MAN_NAME\" run. Deleting it lets Word merge the two runs\n#    into one (and the OOXML writer renumbers the remaining bookmark\n#    ids automatically).\n$bm = $d.Bookmarks.Item(\"__DdeLink__72_1572625028\")\n$bm.Delete()\n\n# 2) Strip the surrounding quote characters around $MAN_NAME so the\n#    line reads \"\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: $MAN_NAME\" instead of '\u0412\u0438\u0440\u043e\u0431\u043d\u0438\u043a: \"$MAN_NAME\"'.\n#    (This also merges the two remaining runs into one, since Word\n#    re-runs the replaced span as a single run with the original\n#    run's formatting.)\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '\"$MAN_NAME\"'\n$find.Replacement.Text = '$MAN_NAME'\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute(\n    $find.Text,        # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    $wdFindContinue,   # Wrap\n    $false,            # Format\n    $find.Replacement.Text,  # ReplaceWith\n    $wdReplaceAll      # Replace\n) | Out-Null\n"}
